$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting (border/bold/center/top alignment) from A10 onto the new
# A11 label cell, then set its text. New shared string "Q9" is added
# implicitly when the value is assigned.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Q9"

# Update rows 2-10 (shift of one row: old data for Q7..Q14 replaced by new
# forecast-error values; row indices correspond to sheet rows unchanged).
$data = @(
    @(2,  0.08955173966697957, 2.273904322189158, 8.17078559450975,  2.858458604652121,  2.885484604367352, 51),
    @(3,  0.4463754146773374,  1.740189769409006, 5.865742303203072, 2.421929458758671,  2.404606822707389, 50),
    @(4,  0.2019302632364671,  1.806917761578944, 8.086275077752271, 2.843637648813975,  2.865853069188615, 49),
    @(5,  0.4252209318853117,  2.131188842423581, 8.399318469451574, 2.89815777166316,   2.897130899773421, 48),
    @(6,  0.2652521615722713,  2.002825916836749, 5.954243375372193, 2.440131835654007,  2.45189624624986,  47),
    @(7,  0.4999690279088892,  2.08768557276191,  7.593124772848787, 2.755562514777842,  2.739769512702802, 46),
    @(8,  0.2296267861977144,  2.088810190185991, 8.927205418845755, 2.987842937445969,  3.012668188838675, 45),
    @(9,  0.6094138230921922,  1.868813072520197, 6.152424372330578, 2.480408106004046,  2.432176378436287, 44),
    @(10, 0.346472456232026,   2.156096895354216, 7.27551841264675,  2.697316891402779,  2.706629584565636, 43),
    @(11, 0.7689217661775961,  2.340507979975441, 8.372201904652238, 2.893475748067061,  2.823250035713138, 42)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
